$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3608   # was 3596
$ws.Range("F5").Value = 3608   # was 3596
$ws.Range("F6").Value = 265   # was 262
$ws.Range("F7").Value = 5139   # was 5129
$ws.Range("F8").Value = 540   # was 536
$ws.Range("F9").Value = 372   # was 365
$ws.Range("F10").Value = 201   # was 200
$ws.Range("F11").Value = 697   # was 695
$ws.Range("F13").Value = 98   # was 97
$ws.Range("F14").Value = 36   # was 35
$ws.Range("F15").Value = 708   # was 705
$ws.Range("F16").Value = 321   # was 320
$ws.Range("F21").Value = 363   # was 362
$ws.Range("F22").Value = 4932   # was 4925
$ws.Range("F26").Value = 6059   # was 6057
$ws.Range("F29").Value = 3226   # was 3225
$ws.Range("F31").Value = 716   # was 715
$ws.Range("F34").Value = 125   # was 123
$ws.Range("F35").Value = 141   # was 140
$ws.Range("F36").Value = 1042   # was 1034
$ws.Range("F37").Value = 83   # was 82
$ws.Range("F38").Value = 25   # was 24
$ws.Range("F40").Value = 877   # was 875
$ws.Range("F41").Value = 1021   # was 1014
$ws.Range("F42").Value = 2033   # was 2031

# --- Sheet: 演出 (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 47   # was 46
$ws.Range("F5").Value = 56   # was 55

# --- Sheet: 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 3608   # was 3596
$ws.Range("F8").Value = 3608   # was 3596
$ws.Range("F9").Value = 265   # was 262
$ws.Range("F10").Value = 5139   # was 5129
$ws.Range("F11").Value = 540   # was 536
$ws.Range("F12").Value = 372   # was 365
$ws.Range("F13").Value = 201   # was 200
$ws.Range("F14").Value = 697   # was 695
$ws.Range("F16").Value = 98   # was 97
$ws.Range("F17").Value = 36   # was 35
$ws.Range("F18").Value = 708   # was 705
$ws.Range("F19").Value = 321   # was 320
$ws.Range("F21").Value = 47   # was 46
$ws.Range("F25").Value = 363   # was 362
$ws.Range("F26").Value = 4932   # was 4925
$ws.Range("F30").Value = 6059   # was 6057
$ws.Range("F33").Value = 3226   # was 3225
$ws.Range("F35").Value = 716   # was 715
$ws.Range("F39").Value = 125   # was 123
$ws.Range("F40").Value = 141   # was 140
$ws.Range("F41").Value = 1042   # was 1034
$ws.Range("F42").Value = 83   # was 82
$ws.Range("F43").Value = 25   # was 24
$ws.Range("F45").Value = 877   # was 875
$ws.Range("F46").Value = 1021   # was 1014
$ws.Range("F48").Value = 2033   # was 2031
$ws.Range("F50").Value = 56   # was 55

$wb.Save()
